{"js": "// Add the three new character styles used by the edit (GaNStyle, GaNParagraph, GaNLinks).\ncontext.document.addStyle(\"GaNStyle\", Word.StyleType.character);\ncontext.document.addStyle(\"GaNParagraph\", Word.StyleType.character);\ncontext.document.addStyle(\"GaNLinks\", Word.StyleType.character);\nawait context.sync();\n\nconst styles = context.document.getStyles();\n\nconst gaNStyle = styles.getByNameOrNullObject(\"GaNStyle\");\ngaNStyle.font.name = \"Calibri\";\ngaNStyle.font.size = 14;\n\nconst gaNParagraph = styles.getByNameOrNullObject(\"GaNParagraph\");\ngaNParagraph.font.name = \"Calibri\";\ngaNParagraph.font.size = 10;\n\nconst gaNLinks = styles.getByNameOrNullObject(\"GaNLinks\");\ngaNLinks.font.name = \"Calibri\";\ngaNLinks.font.bold = true;\ngaNLinks.font.color = \"#000080\";\ngaNLinks.font.size = 9.5;\ngaNLinks.font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// The \"Informace v t\u00e9to p\u0159\u00edru\u010dce...\" paragraph appears 4 times in the body; each one\n// repeats a stray trailing fragment (\"16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna\")\n// that must be dropped, and the run restyled with the new GaNStyle character style.\nconst body = context.document.body;\nconst oldText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, \" +\n  \"14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00edSouhv\u011bzd\u00ed \" +\n  \"Orion.16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna\";\nconst newText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, \" +\n  \"14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00edSouhv\u011bzd\u00ed \" +\n  \"Orion.\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const r of results.items) {\n  r.insertText(newText, \"Replace\");\n  await context.sync();\n  r.style = \"GaNStyle\";\n  await context.sync();\n}\n\n// The \"Jen\u00edk Hollan, CzechGlobe (\u2026)\" credit run gets the GaNLinks character style.\nconst linkResults = body.search(\n  \"Jen\u00edk Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/\",\n  { matchCase: true }\n);\nlinkResults.load(\"items\");\nawait context.sync();\n\nfor (const r of linkResults.items) {\n  r.style = \"GaNLinks\";\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Add the three new character styles used by the edit (GaNStyle, GaNParagraph, GaNLinks).\n# wdStyleTypeCharacter = 2\n$s1 = $d.Styles.Add(\"GaNStyle\", 2)\n$s1.Font.Name = \"Calibri\"\n$s1.Font.Size = 14\n\n$s2 = $d.Styles.Add(\"GaNParagraph\", 2)\n$s2.Font.Name = \"Calibri\"\n$s2.Font.Size = 10\n\n$s3 = $d.Styles.Add(\"GaNLinks\", 2)\n$s3.Font.Name = \"Calibri\"\n$s3.Font.Bold = $true\n$s3.Font.Color = 8388608 # RGB 000080 (navy) stored BGR in the OLE_COLOR\n$s3.Font.Size = 9.5\n$s3.Font.Underline = 1 # wdUnderlineSingle\n\n# The \"Informace v t\u00e9to p\u0159\u00edru\u010dce...\" paragraph appears 4 times in the body; each one\n# repeats a stray trailing fragment (\"16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna\")\n# that must be dropped, and the run restyled with the new GaNStyle character style.\n$oldText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00edSouhv\u011bzd\u00ed Orion.16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna\"\n$newText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16. \u2013 25. ledna, 14. \u2013 23. \u00fanora, 14. \u2013 24. b\u0159ezna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00edSouhv\u011bzd\u00ed Orion.\"\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text -like \"*Informace v t\u00e9to p\u0159\u00edru\u010dce*\") {\n        # Exclude the trailing paragraph mark so the style lands on the run (rStyle),\n        # not the paragraph (pStyle).\n        $start = $r.Start\n        $textRange = $d.Range($start, $r.End - 1)\n        $textRange.Text = $newText\n        $styledRange = $d.Range($start, $start + $newText.Length)\n        $styledRange.Style = \"GaNStyle\"\n    }\n}\n\n# The \"Jen\u00edk Hollan, CzechGlobe (\u2026)\" credit run gets the GaNLinks character style.\n$linkRange = $d.Content\n$found = $linkRange.Find.Execute(\"Jen\u00edk Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/\")\nif ($found) {\n    $linkRange.Style = \"GaNLinks\"\n}\n"}
